# ============================================================================
# Additional scraping: add a "Player Info" sheet, convert the
# MATCH_CARD_LINK url columns on the existing batting/bowling sheets into a
# plain MATCH_CODE number-as-text column, drop the now-pointless empty
# INNING_NUMBER placeholder cells, and append a new "ODI Batting Extra"
# sheet with the extra per-innings batting detail that was scraped.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# Helper: turn a header/data range into the bold + thin-border + centered
# style used by every header row in this workbook.
# ----------------------------------------------------------------------------
function Set-HeaderStyle($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108   # xlCenter
    $range.VerticalAlignment = -4160     # xlTop
    $range.Borders.LineStyle = 1         # xlContinuous
}

# ----------------------------------------------------------------------------
# IMPORTANT: this COM host anchors worksheet handles by POSITION. Inserting
# a sheet before a position a handle is anchored at silently re-seats that
# handle onto the newly inserted sheet. So: do ALL Worksheets.Add() calls
# first (in an order that never inserts before an already-created handle),
# and only fetch+use sheet handles (fresh, by name) afterwards.
# ----------------------------------------------------------------------------

# 1) Insert "Player Info" right before "ODI Batting" (i.e. new first sheet).
$battingSheetTmp = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheetTmp)
$playerInfo.Name = "Player Info"

# 2) Append "ODI Batting Extra" after the current last sheet ("ODI Bowling").
#    This happens strictly after every other sheet's position, so it does
#    not disturb $playerInfo (position 1).
$lastSheetTmp = $wb.Worksheets.Item($wb.Worksheets.Count)
$extraSheet = $wb.Worksheets.Add($null, $lastSheetTmp)
$extraSheet.Name = "ODI Batting Extra"

# From here on, always re-fetch sheets fresh by name before using them.

# ----------------------------------------------------------------------------
# Populate "Player Info"
# ----------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Item("Player Info")

$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $piHeaders.Length; $c++) {
    $cell = $playerInfo.Cells.Item(1, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $piHeaders[$c - 1]
}
Set-HeaderStyle($playerInfo.Range("A1:D1"))

$piRow = @("4956", "Hayden Rashidi Walsh", "Left Handed", "Right Arm Leg Break")
for ($c = 1; $c -le $piRow.Length; $c++) {
    $cell = $playerInfo.Cells.Item(2, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $piRow[$c - 1]
}

# ----------------------------------------------------------------------------
# "ODI Batting": MATCH_CARD_LINK (col D) -> MATCH_CODE, and drop the
# empty INNING_NUMBER (col B) placeholder cells for "did not bat" rows.
# ----------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Cells.Item(1, 4).Value = "MATCH_CODE"

$battingLastRow = $battingSheet.UsedRange.Rows.Count
for ($r = 2; $r -le $battingLastRow; $r++) {
    $linkCell = $battingSheet.Cells.Item($r, 4)
    $txt = $linkCell.Text
    if ($txt -and $txt.Contains("MatchCode=")) {
        $code = $txt.Substring($txt.LastIndexOf("=") + 1)
        $linkCell.NumberFormat = "@"
        $linkCell.Value = $code
    }

    $inningCell = $battingSheet.Cells.Item($r, 2)
    $inningTxt = $inningCell.Text
    if (-not $inningTxt) {
        $inningCell.ClearContents()
    }
}

# ----------------------------------------------------------------------------
# "ODI Bowling": MATCH_CARD_LINK (col B) -> MATCH_CODE.
# ----------------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Cells.Item(1, 2).Value = "MATCH_CODE"

$bowlingLastRow = $bowlingSheet.UsedRange.Rows.Count
for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $linkCell = $bowlingSheet.Cells.Item($r, 2)
    $txt = $linkCell.Text
    if ($txt -and $txt.Contains("MatchCode=")) {
        $code = $txt.Substring($txt.LastIndexOf("=") + 1)
        $linkCell.NumberFormat = "@"
        $linkCell.Value = $code
    }
}

# ----------------------------------------------------------------------------
# Populate "ODI Batting Extra"
# ----------------------------------------------------------------------------
$extraSheet = $wb.Worksheets.Item("ODI Batting Extra")

$exHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $exHeaders.Length; $c++) {
    $cell = $extraSheet.Cells.Item(1, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $exHeaders[$c - 1]
}
Set-HeaderStyle($extraSheet.Range("A1:F1"))

# MATCH_CODE, BATTING_POSITION (numeric, blank = $null), NUM_4, NUM_6,
# PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH ($null = blank placeholder cell)
$exData = @(
    @("4378", 10, "1", "0", "1.62%", "NO"),
    @("4379", $null, $null, $null, $null, "NO"),
    @("4385", 9, "4", "0", "19.01%", "NO"),
    @("4391", 9, $null, $null, $null, "NO"),
    @("4394", 9, "2", "1", "6.92%", "NO"),
    @("4397", 8, "0", "0", "0.66%", "NO"),
    @("4413", 9, "4", "0", "16.26%", "NO"),
    @("4417", 9, $null, $null, $null, "NO"),
    @("4483", $null, $null, $null, $null, "NO"),
    @("4484", 10, "1", "0", "7.69%", "NO"),
    @("4486", 10, $null, $null, $null, "NO"),
    @("4536", 10, $null, $null, $null, "NO"),
    @("4577", 9, $null, $null, $null, "NO"),
    @("4580", $null, $null, $null, $null, "NO"),
    @("4583", $null, $null, $null, $null, "NO"),
    @("4586", 10, "0", "0", "1.39%", "NO"),
    @("4590", $null, $null, $null, $null, "NO"),
    @("4592", 10, "2", "0", "7.30%", "NO"),
    @("4623", $null, $null, $null, $null, $null),
    @("4624", $null, $null, $null, $null, $null)
)

for ($i = 0; $i -lt $exData.Length; $i++) {
    $r = $i + 2
    $rowData = $exData[$i]

    $codeCell = $extraSheet.Cells.Item($r, 1)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $rowData[0]

    $posVal = $rowData[1]
    if ($null -ne $posVal) {
        $extraSheet.Cells.Item($r, 2).Value = $posVal
    }

    for ($c = 3; $c -le 5; $c++) {
        $v = $rowData[$c - 1]
        if ($null -ne $v) {
            $cell = $extraSheet.Cells.Item($r, $c)
            $cell.NumberFormat = "@"
            $cell.Value = $v
        }
    }

    $momVal = $rowData[5]
    if ($null -ne $momVal) {
        $momCell = $extraSheet.Cells.Item($r, 6)
        $momCell.NumberFormat = "@"
        $momCell.Value = $momVal
    }
}

Write-Output "done"
